$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E9"  = 27
    "E15" = 162
    "E16" = 12
    "F17" = 61
    "H17" = 93
    "E18" = 114
    "F18" = 51
    "H18" = 87
    "E19" = 58
    "E28" = 19
    "F28" = 14
    "H28" = 16
    "E34" = 22
    "E40" = 21
    "E41" = 41
    "F41" = 17
    "H41" = 28
    "E46" = 30
    "E47" = 61
    "F47" = 38
    "H47" = 49
    "E49" = 74
    "E56" = 8
    "F56" = 3
    "H56" = 5
    "E61" = 29
    "E62" = 46
    "E63" = 38
    "E65" = 33
    "E67" = 40
    "F67" = 24
    "H67" = 33
    "E71" = 38
    "F71" = 18
    "H71" = 28
    "E72" = 45
    "F72" = 24
    "H72" = 35
    "E73" = 29
    "E76" = 52
    "F76" = 19
    "H76" = 36
    "E83" = 11
    "E88" = 25
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
